$d = $word.ActiveDocument

# Locate the exact split point: right before "and greater park ecosystems across"
# within "...between parks and greater park ecosystems across canada..."
$rng = $d.Content
$found = $rng.Find.Execute("between parks and greater park ecosystems across", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$splitOffset = $rng.Start + "between parks ".Length
$insertPoint = $d.Range($splitOffset, $splitOffset)

$d.Bookmarks.Add("_GoBack", $insertPoint)
